$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")
$ws.Activate()

# Add a new row of data describing the "esd" dataset (English and Mandarin),
# specifying iso 3166-1 locale codes alongside the existing entries.
$ws.Range("A16").Value = "esd"
$ws.Range("B16").Value = "acted"
$ws.Range("H16").Value = "English and Mandarin"

# Leave selection where Excel would naturally land after entering the row.
$ws.Range("B17").Select() | Out-Null
